# Weekly update to "Hortaliza, Vega Modelo de Temuco - Papa" sheet.
# Inserts two new daily-price rows (686 and 687) above the existing data,
# which pushes the previously-existing rows 686-735 down to 688-737.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 686 (shifts old rows 686:735 down to 688:737)
$ws.Range("A686:A687").EntireRow.Insert()

# New row 686: Patagonia, 1a (guarda)
$ws.Cells.Item(686, 1).Value = 10
$ws.Cells.Item(686, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(686, 3).Value = "La Araucanía"
$ws.Cells.Item(686, 4).Value = 44706
$ws.Cells.Item(686, 5).Value = 9
$ws.Cells.Item(686, 6).Value = 100114001
$ws.Cells.Item(686, 7).Value = "Papa"
$ws.Cells.Item(686, 8).Value = "Patagonia"
$ws.Cells.Item(686, 9).Value = "1a (guarda)"
$ws.Cells.Item(686, 10).Value = 500
$ws.Cells.Item(686, 11).Value = 8000
$ws.Cells.Item(686, 12).Value = 8000
$ws.Cells.Item(686, 13).Value = 8000
$ws.Cells.Item(686, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(686, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(686, 16).Value = 320
$ws.Cells.Item(686, 17).Value = 25
$ws.Cells.Item(686, 18).Value = "Hortaliza"

# New row 687: Rosara, 1a (guarda)
$ws.Cells.Item(687, 1).Value = 10
$ws.Cells.Item(687, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(687, 3).Value = "La Araucanía"
$ws.Cells.Item(687, 4).Value = 44706
$ws.Cells.Item(687, 5).Value = 9
$ws.Cells.Item(687, 6).Value = 100114001
$ws.Cells.Item(687, 7).Value = "Papa"
$ws.Cells.Item(687, 8).Value = "Rosara"
$ws.Cells.Item(687, 9).Value = "1a (guarda)"
$ws.Cells.Item(687, 10).Value = 400
$ws.Cells.Item(687, 11).Value = 8000
$ws.Cells.Item(687, 12).Value = 8000
$ws.Cells.Item(687, 13).Value = 8000
$ws.Cells.Item(687, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(687, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(687, 16).Value = 320
$ws.Cells.Item(687, 17).Value = 25
$ws.Cells.Item(687, 18).Value = "Hortaliza"
